# Adjusted card layout for large description and two cards added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (existing "Relance de commandement" card): description + timing text replaced ---
# D2: full description text (now HTML-entity-escaped-as-plain-text, matches new source content)
$ws.Range("D2").Value = "Utilisez ce stratag&egrave;me apr&egrave;s que vous ayez effectu&eacute; un jet de touche, blessure, d&eacute;gats, sauvegarde, Avance, charge, psy, abjurer ou nombre d'attaque effectu&eacute;es par une arme. <strong>Relancez ce jet ou ce test</strong>."
# B2: short timing label text
$ws.Range("B2").Value = "Apr&egrave; jet touche, svg, avance, charge, psy, abjurer, nombre attaque"

# Row 2 grew taller to fit the longer description
$ws.Rows.Item(2).RowHeight = 75

# --- Row 3 (new card: "Exterminez-les") ---
$ws.Range("A3").Value = "Exterminez-les"
$ws.Range("B3").Value = "Unit&eacute; ennemie bat en retraite"
$ws.Range("C3").Value = "strats-sources\assets\timer.png"
$ws.Range("D3").Value = "Utilisez ce stratag&egrave;me quand une unit&eacute; ennemie bat en retraite, avant que la moindre figurine de l'unt&eacute; soit d&eacute;plac&eacute;e.<br/>Jetez 1D6 pour chaque figurine se trouvant &agrave; d'engagement de l'unit&eacute; ennemie ; pour chaque r&eacute;sultat de 6, l'unit&eacute; ennemie subit 1 blessure mortelle."
$ws.Range("E3").Value = "255"
$ws.Range("F3").Value = "1"
$ws.Range("G3").Value = "Base"
$ws.Range("H3").Value = "strats-sources\assets\base-rulebook.png"
$ws.Rows.Item(3).RowHeight = 90

# --- Row 4 (new card: "Repli desespere") ---
$ws.Range("A4").Value = "Repli desesp&eacute;r&eacute;"
$ws.Range("B4").Value = "Phase de mouvement"
$ws.Range("C4").Value = "strats-sources\assets\timer.png"
$ws.Range("E4").Value = "255"
$ws.Range("F4").Value = "2"
$ws.Range("D4").Value = "Choisir une unit&eacute; qui ne s'est pas d&eacute;plac&eacute; &agrave; port&eacute;e d'engagement d'aun moins une unit&eacute; ennemie.<br/><strong>Jetez 1D6</strong> pour chaque fig de l'unit&eacute; finie. <strong>Sur 1 l'unit&eacute est d&eacute;truite</strong>. Si l'unit&eacute; n'est pas d&eacute;truite elle peut tenter de <i>Battre en Retraite</i>.<br/>Si la figurine termine son mouvement &agrave; <i>Port&eacute;e d'Engagement</i>, elle d&eacute;truite.<br/><strong>L'unit&eacute; ne peut plus rien faire</strong>, m&ecirc;me si elle a une r&egrave;gle qui lui permet de faire actions apr&egrave; avoir Battu en Retraite."
$ws.Range("G4").Value = "Base"
$ws.Range("H4").Value = "strats-sources\assets\base-rulebook.png"
$ws.Rows.Item(4).RowHeight = 165

# --- Scroll / selection state ---
$ws.Range("D5").Select() | Out-Null
